# Word COM-interop script applying the "Final Checking Part I" edits to the
# family court Occupation Order document.
#
# NOTE on ordering: calling Table.Rows.Add() on this runtime invalidates
# $d.Paragraphs index-based lookups done afterwards, so all paragraph-level
# work (Find/Replace text edits + paragraph insertions) is done first, and
# the two new child-table rows are appended at the very end.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
}

# 1. Heading: Family Law Act year
Replace-Text "Family Law Act 1997" "Family Law Act 1996"

# 2. Children table header row text
Replace-Text "The full name of the child" "The full names of the children"
Replace-Text "Date of Birth" "Dates of Birth"

# 3. Important notice box
Replace-Text "IMPORTANT NOTICE TO THE RESPONDENT     of ,      ." "IMPORTANT NOTICE TO THE RESPONDENT JOSEPH   SAMUELSON  of 2 THE DRIVE,    CROYDON  CR218 4TN."
Replace-Text "WARNING: ALTERNATIVELY, IF YOU DISOBEY" "WARNING: IF YOU DISOBEY"

# 4. Before District Judge line
Replace-Text "Before District Judge ." "Before Deputy District Judge Campbell  in private on  13 February 2025  at an interim hearing ."

# 5. Parties line
Replace-Text "The applicant is   represented by 0   The respondent is   represented by 0  " "The applicant is Samantha Samuels represented by  The respondent is Joseph  Samuelson represented by "

# 6. Relevant child intro line
Replace-Text "The relevant child within the meaning of the Family Law Act 1995 is:" "The `"relevant children`" within the meaning of Family Law Act 1996 are:"

# 7. Item (a) child details
Replace-Text "  a  born on 06 Feb 2015" " a boy born on 01 Jan 2021"

# 8. Insert items (b) and (c) right after item (a) -- paragraph 40 in the
#    original document; none of the edits above add/remove paragraphs so the
#    index is still valid here.
$pA = $d.Paragraphs.Item(40)
$pA.Range.InsertParagraphAfter()
$pB = $d.Paragraphs.Item(41)
$pB.Range.Text = "`tb`t a girl born on 02 Feb 2022"
$pB.Range.InsertParagraphAfter()
$pC = $d.Paragraphs.Item(42)
$pC.Range.Text = "`tc`t a boy born on 03 Mar 2020"

# 9. Family home line
Replace-Text "The `"family home`" is the property at    " "The `"family home`" is the property at The larches East Hampton Road Hailsham BN28 480"

# 10. Recitals - occupation order made against
Replace-Text "This is an  Occupation Order made against the respondent    on  by District Judge  on the application of the applicant   " "This is an  Occupation Order made against the respondent Joseph   Samuelson on 13 February 2025 by Deputy District Judge Campbell on the application of the applicant Samantha  Samuels"

# 11. IT IS ORDERED declarations
Replace-Text "The court declares that the applicant , has home rights in , , , ." "The court declares that the applicant Samantha Samuels, has home rights in The larches, East Hampton Road, Hailsham, BN28 480."

Replace-Text "The court declares that the  's home rights in , , ,  shall not end when the respondent   dies or their marriage is dissolved and shall continue until the determination of the applicant's financial provision claims or a further order is made." "The court declares that the Samantha Samuels's home rights in The larches, East Hampton Road, Hailsham, BN28 480 shall not end when the respondent Joseph  Samuelson dies or their marriage is dissolved and shall continue until the determination of the applicant's financial provision claims or a further order is made."

Replace-Text "The respondent   must not use or threaten violence against the applicant   and must not instruct, encourage or in any way suggest any other person should do so." "The respondent Joseph  Samuelson must not use or threaten violence against the applicant Samantha Samuels and must not instruct, encourage or in any way suggest any other person should do so."

Replace-Text "The respondent   must not intimidate, harass or pester the applicant    and must not instruct, encourage or in any way suggest any other person should do so." "The respondent Joseph  Samuelson must not intimidate, harass or pester the applicant Samantha Samuels  and must not instruct, encourage or in any way suggest any other person should do so."

Replace-Text "The respondent,  , must not telephone, text, email or otherwise contact the applicant  , including via social networking websites or other forms of electronic messaging." "The respondent, Joseph  Samuelson, must not telephone, text, email or otherwise contact the applicant Samantha Samuels, including via social networking websites or other forms of electronic messaging."

# 12. The next two paragraphs get entirely new content (damage-to-property
#     undertakings) -- paragraphs 53 and 54 in the original numbering, now
#     shifted to 55 and 56 because items (b)/(c) added two paragraphs above.
$p55 = $d.Paragraphs.Item(55)
$p55.Range.Text = "The respondent, Joseph  Samuelson, must not damage, attempt to damage or threaten to damageany property owned by or in the possession or controlof the applicant Samantha Samuels,and must not instruct,encourage or in any way suggest that any other person should do so."

$p56 = $d.Paragraphs.Item(56)
$p56.Range.Text = "The respondent Joseph  Samuelson, must not damage, attempt to damage or threaten to damagethe property or contents of The larches, East Hampton Road, Hailsham, BN28 480 and must not instruct, encourage or in any way suggestthat any other person should do so."

# 13. Insert six new ListNumber paragraphs after paragraph 56
$newParaTexts = @(
    "The respondent Joseph  Samuelson, must not go to, enter or attempt to enter The larches, East Hampton Road, Hailsham, BN28 480 or  any property where he believes the applicant Samantha Samuels to be living.",
    "The Respondent Joseph  Samuelson, must not use or threaten violence against the relevant children  and must not instruct, encourage or in any way suggest that any other person should do so.",
    "The respondent Joseph  Samuelson, must not telephone, text,email or otherwise contact or attempt  the relevant children including via social networking websites or other forms of electronic messaging.",
    "The respondent Joseph  Samuelson, must not go to, enter or attempt to enter the school premises known as Hellingly Primary School, Hellingly, East Sussex BN27 1PQ except by prior written invitation from the school authorities.",
    "This order shall be effective against the respondent Joseph  Samuelson once it is personally served on him.",
    "This order shall last until 13 August 2025 unless it is set aside or varied before then by an order of the court."
)

$curIdx = 56
foreach ($t in $newParaTexts) {
    $p = $d.Paragraphs.Item($curIdx)
    $p.Range.InsertParagraphAfter()
    $curIdx = $curIdx + 1
    $newP = $d.Paragraphs.Item($curIdx)
    $newP.Range.Text = $t
}

# 14. Signature block. After the inserts above, the document now has
#     58 (original) + 2 (b/c) + 6 (new undertakings) = 66 paragraphs, with
#     "District Judge " (originally #56) now at #64, and the trailing blank
#     paragraph (originally #58) now at #66.
$pJudge = $d.Paragraphs.Item($curIdx + 2)
$pJudge.Range.Text = "Deputy District Judge Campbell"

$pDate = $d.Paragraphs.Item($d.Paragraphs.Count)
$pDate.Range.Text = "13 February 2025"

# 15. Children table: fill in the existing data row and append two more rows.
#     Done last because Table.Rows.Add() invalidates $d.Paragraphs index
#     lookups afterwards on this runtime.
$childTable = $d.Tables.Item(2)

$row1 = $childTable.Rows.Item(2)
$row1.Cells.Item(1).Range.Text = "Arthur Simpson"
$row1.Cells.Item(2).Range.Text = "boy"
$row1.Cells.Item(3).Range.Text = "01 Jan 2021"

$row2 = $childTable.Rows.Add()
$row2.Cells.Item(1).Range.Text = "Belinda Simpson"
$row2.Cells.Item(2).Range.Text = "girl"
$row2.Cells.Item(3).Range.Text = "02 Feb 2022"

$row3 = $childTable.Rows.Add()
$row3.Cells.Item(1).Range.Text = "Charlie Simpson"
$row3.Cells.Item(2).Range.Text = "boy"
$row3.Cells.Item(3).Range.Text = "03 Mar 2020"

Write-Output "done"
